$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.436.22"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").Value = "3.566.76"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "620.96"
$ws.Range("E5").Value = "  +2.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.07"
$ws.Range("E6").Value = "  +3.65%  "

$ws.Range("D7").Value = "3.565.16"
$ws.Range("E7").Value = "  +1.82%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  +2.28%  "

$ws.Range("E10").Value = "  +5.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.39"
$ws.Range("E11").Value = "  +5.85%  "

$ws.Range("E12").Value = "  +3.91%  "

$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.23"
$ws.Range("E14").Value = "  +5.29%  "

$ws.Range("D15").Value = "4.169.51"
$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("D16").Value = "3.565.62"
$ws.Range("E16").Value = "  +1.76%  "

$ws.Range("D17").Value = "68.428.74"
$ws.Range("E17").Value = "  +1.63%  "

$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("E19").Value = "  +6.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.03"
$ws.Range("E20").Value = "  +6.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("E21").Value = "  +11.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.27"
$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.643"
$ws.Range("E23").Value = "  +3.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.55"
$ws.Range("E24").Value = "  +1.46%  "

$ws.Range("E25").Value = "  +1.31%  "

$ws.Range("D26").Value = "3.708.88"
$ws.Range("E26").Value = "  +1.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.19%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.19"
$ws.Range("E28").Value = "  +11.21%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.55"
$ws.Range("E29").Value = "  +3.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.71"
$ws.Range("E30").Value = "  +10.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.56"
$ws.Range("E31").Value = "  +3.15%  "

$ws.Range("E32").Value = "  +3.79%  "

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.36"
$ws.Range("E34").Value = "  +3.27%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.17"
$ws.Range("E35").Value = "  +1.75%  "

$ws.Range("E36").Value = "  +3.49%  "

$ws.Range("D37").Value = "3.559.77"
$ws.Range("E37").Value = "  +1.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.25"

$ws.Range("E39").Value = "  +8.76%  "

$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "180.85"
$ws.Range("E41").Value = "  +4.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0917"
$ws.Range("E42").Value = "  +4.74%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.59"
$ws.Range("E44").Value = "  +3.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.13"
$ws.Range("E45").Value = "  +14.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.899"
$ws.Range("E46").Value = "  +2.07%  "

$ws.Range("E47").Value = "  +1.69%  "

$ws.Range("E48").Value = "  +4.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.68"
$ws.Range("E49").Value = "  +4.54%  "

$ws.Range("E50").Value = "  +3.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.263"
$ws.Range("E51").Value = "  +7.57%  "
